# Add a 'Replace Existing' column to the sequence upload template, and
# add a sample data row demonstrating its use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D (shifts old D..I -> E..J), which carries over
# the existing column width/style formatting from column C/E.
$ws.Columns("D:D").Insert()

# New header cell for the inserted column.
$ws.Range("D1").Value = "Replace Existing"

# New sample data row (row 2), filled in an order that mirrors how the
# template was originally authored.
$ws.Range("D2").Value = "1"
$ws.Range("F2").Value = "FirstGrowiReady"
$ws.Range("C2").Value = "987659"
$ws.Range("B2").Value = "Lau"
$ws.Range("A2").Value = "Brigette"

# Move the active selection, matching the saved view state.
$ws.Range("G4").Select()

Write-Output "Added 'Replace Existing' column and sample row."
